# Update cryptocurrency price/volume data per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D, E
$ws.Range("D2").Value = '90.540.31'
$ws.Range("E2").Value = '  +5.73%  '

# Row 3: update D, E
$ws.Range("D3").Value = '3.270.10'
$ws.Range("E3").Value = '  +0.56%  '

# Row 4: update E
$ws.Range("E4").Value = '  -0.07%  '

# Row 5: update D, E
$ws.Range("D5").Value = '''214.80'
$ws.Range("E5").Value = '  +2.94%  '

# Row 6: update D, E
$ws.Range("D6").Value = '''625.71'
$ws.Range("E6").Value = '  +0.52%  '

# Row 7: update D, E
$ws.Range("D7").Value = '''0.412'
$ws.Range("E7").Value = '  +15.19%  '

# Row 8: update D, E
$ws.Range("D8").Value = '''0.710'
$ws.Range("E8").Value = '  +9.41%  '

# Row 9: update D, E
$ws.Range("D9").Value = '''0.998'
$ws.Range("E9").Value = '  -0.26%  '

# Row 10: update D, E
$ws.Range("D10").Value = '3.268.43'
$ws.Range("E10").Value = '  +1.00%  '

# Row 11: update D, E
$ws.Range("D11").Value = '''0.593'
$ws.Range("E11").Value = '  +4.15%  '

# Row 12: update D, E
$ws.Range("D12").Value = '''0.0000276'
$ws.Range("E12").Value = '  +9.84%  '

# Row 13: update E
$ws.Range("E13").Value = '  +2.15%  '

# Row 14: update B, C, D, E
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''34.20'
$ws.Range("E14").Value = '  +1.81%  '

# Row 15: update B, C, D, E
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = '''5.44'
$ws.Range("E15").Value = '  +3.90%  '

# Row 16: update D, E
$ws.Range("D16").Value = '3.841.27'
$ws.Range("E16").Value = '  -0.28%  '

# Row 17: update D, E
$ws.Range("D17").Value = '90.023.16'
$ws.Range("E17").Value = '  +5.18%  '

# Row 18: update D, E
$ws.Range("D18").Value = '3.261.84'
$ws.Range("E18").Value = '  +0.31%  '

# Row 19: update D, E
$ws.Range("D19").Value = '''3.26'
$ws.Range("E19").Value = '  +10.14%  '

# Row 20: update D, E
$ws.Range("D20").Value = '''14.18'
$ws.Range("E20").Value = '  +1.96%  '

# Row 21: update D, E
$ws.Range("D21").Value = '''432.28'
$ws.Range("E21").Value = '  +1.74%  '

# Row 22: update D, E
$ws.Range("D22").Value = '''8.95'
$ws.Range("E22").Value = '  +1.02%  '

# Row 23: update B, C, D, E
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").Value = '''0.0000191'
$ws.Range("E23").Value = '  +50.96%  '

# Row 24: update B, C, D, E
$ws.Range("B24").Value = 'Polkadot'
$ws.Range("C24").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D24").Value = '''5.30'
$ws.Range("E24").Value = '  +0.80%  '

# Row 25: update D, E
$ws.Range("D25").Value = '''5.43'
$ws.Range("E25").Value = '  +7.98%  '

# Row 26: update D, E
$ws.Range("D26").Value = '''12.04'
$ws.Range("E26").Value = '  -2.04%  '

# Row 27: update D, E
$ws.Range("D27").Value = '3.398.49'
$ws.Range("E27").Value = '  -0.99%  '

# Row 28: update D, E
$ws.Range("D28").Value = '''76.00'
$ws.Range("E28").Value = '  +1.40%  '

# Row 29: update E
$ws.Range("E29").Value = '  -0.03%  '

# Row 30: update D, E
$ws.Range("D30").Value = '''0.175'
$ws.Range("E30").Value = '  +2.77%  '

# Row 31: update D, E
$ws.Range("D31").Value = '''0.996'
$ws.Range("E31").Value = '  -0.71%  '

# Row 32: update D, E
$ws.Range("D32").Value = '''570.16'
$ws.Range("E32").Value = '  +5.35%  '

# Row 33: update D, E
$ws.Range("D33").Value = '''8.61'
$ws.Range("E33").Value = '  -1.19%  '

# Row 34: update D, E
$ws.Range("D34").Value = '''7.19'
$ws.Range("E34").Value = '  +6.21%  '

# Row 35: update D, E
$ws.Range("D35").Value = '''1.36'
$ws.Range("E35").Value = '  -3.42%  '

# Row 36: update D, E
$ws.Range("D36").Value = '''1.91'
$ws.Range("E36").Value = '  -1.29%  '

# Row 37: update D, E
$ws.Range("D37").Value = '''3.54'
$ws.Range("E37").Value = '  +21.97%  '

# Row 38: update B, C, D, E
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '''22.64'
$ws.Range("E38").Value = '  +1.96%  '

# Row 39: update B, C, D, E
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.134'
$ws.Range("E39").Value = '  -1.12%  '

# Row 40: update D, E
$ws.Range("D40").Value = '''22.31'
$ws.Range("E40").Value = '  +3.41%  '

# Row 41: update D, E
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.58%  '

# Row 42: update D, E
$ws.Range("D42").Value = '''0.394'
$ws.Range("E42").Value = '  +1.72%  '

# Row 43: update D, E
$ws.Range("D43").Value = '''1.99'
$ws.Range("E43").Value = '  +1.45%  '

# Row 44: update E
$ws.Range("E44").Value = '  +0.09%  '

# Row 45: update D, E
$ws.Range("D45").Value = '''151.08'
$ws.Range("E45").Value = '  -4.37%  '

# Row 46: update D, E
$ws.Range("D46").Value = '''182.75'
$ws.Range("E46").Value = '  +3.42%  '

# Row 47: update D, E
$ws.Range("D47").Value = '''43.93'
$ws.Range("E47").Value = '  +0.09%  '

# Row 48: update D, E
$ws.Range("D48").Value = '''0.129'
$ws.Range("E48").Value = '  +9.43%  '

# Row 49: update D, E
$ws.Range("D49").Value = '''1.29'
$ws.Range("E49").Value = '  +0.20%  '

# Row 50: update D, E
$ws.Range("D50").Value = '''0.626'
$ws.Range("E50").Value = '  +1.91%  '

# Row 51: update D, E
$ws.Range("D51").Value = '''25.16'
$ws.Range("E51").Value = '  +5.43%  '
